# Update "想去人数" (F column) values across the four sheets of the
# "上海-漫展信息" workbook, reflecting a newer scrape of the source site.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell, newValue) updates.
$updates = @{
    "展览"   = @(
        @{ Cell = "F5";  Value = 1760 },
        @{ Cell = "F7";  Value = 324 },
        @{ Cell = "F8";  Value = 626 },
        @{ Cell = "F15"; Value = 1363 },
        @{ Cell = "F19"; Value = 468 },
        @{ Cell = "F21"; Value = 29 },
        @{ Cell = "F24"; Value = 2323 },
        @{ Cell = "F27"; Value = 4347 },
        @{ Cell = "F32"; Value = 2 },
        @{ Cell = "F34"; Value = 939 }
    )
    "演出"   = @(
        @{ Cell = "F8";  Value = 19 },
        @{ Cell = "F22"; Value = 139 }
    )
    "本地生活" = @(
        @{ Cell = "F3";  Value = 2544 },
        @{ Cell = "F10"; Value = 3004 },
        @{ Cell = "F11"; Value = 522 },
        @{ Cell = "F13"; Value = 245 },
        @{ Cell = "F14"; Value = 267 }
    )
    "全部类型" = @(
        @{ Cell = "F2";  Value = 2544 },
        @{ Cell = "F8";  Value = 3004 },
        @{ Cell = "F10"; Value = 245 },
        @{ Cell = "F11"; Value = 1760 },
        @{ Cell = "F13"; Value = 324 },
        @{ Cell = "F14"; Value = 626 },
        @{ Cell = "F22"; Value = 19 },
        @{ Cell = "F30"; Value = 139 },
        @{ Cell = "F31"; Value = 139 },
        @{ Cell = "F32"; Value = 29 },
        @{ Cell = "F39"; Value = 4347 },
        @{ Cell = "F47"; Value = 2 },
        @{ Cell = "F50"; Value = 939 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
